# Apply the changes described by the commit:
# "update template file with README sheet and version number on template sheet"
#
#  1. Rename the original first sheet ("Sheet1") to "template_v0-00".
#  2. Insert a brand-new "README" sheet right after it, containing four
#     lines of text (with hyperlinks) pointing readers at the GitHub repo,
#     the feedback form and the issue tracker.
#  3. Restore/adjust the selections so that the README tab ends up the
#     active tab (as it would after being freshly added & worked on),
#     while the template sheet keeps a plain A2 selection.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the main template sheet -------------------------------------
$template = $wb.Worksheets.Item(1)
$template.Name = "template_v0-00"

# put the template sheet's selection back to A2 (no more tab-selected there
# once README becomes the active tab)
$template.Activate() | Out-Null
$template.Range("A2").Select() | Out-Null

# --- 2. Insert the new README sheet right after the template ---------------
$readme = $wb.Worksheets.Add($null, $template)
$readme.Name = "README"

# Populate the cells in the exact order needed so the shared-string table
# ends up indexed the same way as the target workbook:
#   74 -> "to find metadata ..."
#   75 -> "to make sure you have to most up-to-date version ..."
#   76 -> "to provide feedback ..."
#   77 -> "or create an issue ..."
$readme.Range("A2").Value = "to find metadata associated with this template (""data dictionary""), go to our repo: https://github.com/atlanticcanadacdc/outsideSubmissionTemplate"
$readme.Range("A1").Value = "to make sure you have to most up-to-date version, download this file directly from our repo: https://github.com/atlanticcanadacdc/outsideSubmissionTemplate"
$readme.Range("A4").Value = "to provide feedback on this template, please submit this Google form"
$readme.Range("A5").Value = "or create an issue on the issue tracker"

# Hyperlinks for each of the four lines
$readme.Hyperlinks.Add($readme.Range("A1"), "https://github.com/atlanticcanadacdc/outsideSubmissionTemplate", "", "to make sure you have to most up to date version, download this file from our repo: https://github.com/atlanticcanadacdc/outsideSubmissionTemplate") | Out-Null
$readme.Hyperlinks.Add($readme.Range("A2"), "https://github.com/atlanticcanadacdc/outsideSubmissionTemplate") | Out-Null
$readme.Hyperlinks.Add($readme.Range("A4"), "https://forms.gle/") | Out-Null
$readme.Hyperlinks.Add($readme.Range("A5"), "https://github.com/atlanticcanadacdc/outsideSubmissionTemplate/issues") | Out-Null

# Match the printed page orientation used on the rest of the workbook
$readme.PageSetup.Orientation = 1

# Final selection on the README sheet -> A5 (matches the saved state)
$readme.Activate() | Out-Null
$readme.Range("A5").Select() | Out-Null
